# Apply the test-data update:
#  - remove the SiteWeatherTest and APIWeatherTest sheets
#  - update Test_Locations with a 4-column table (Place/State/CountryCode/Unit)

$wb = $excel.ActiveWorkbook

# Remove the extra worksheets, keep only Test_Locations
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("SiteWeatherTest").Delete()
$wb.Worksheets.Item("APIWeatherTest").Delete()
$excel.DisplayAlerts = $true

$ws = $wb.Worksheets.Item("Test_Locations")

# Clear the old contents (was A1:C3) before writing the new 4-column table
$ws.Range("A1:C3").ClearContents()

# Header row
$ws.Range("A1").Value = "Place"
$ws.Range("B1").Value = "State"
$ws.Range("C1").Value = "CountryCode"
$ws.Range("D1").Value = "Unit"

# Data rows
$ws.Range("A2").Value = "Pune"
$ws.Range("B2").Value = "Maharashtra"
$ws.Range("C2").Value = "IN"
$ws.Range("D2").Value = "metric"

$ws.Range("A3").Value = "New York"
$ws.Range("B3").Value = "New York"
$ws.Range("C3").Value = "US"
$ws.Range("D3").Value = "metric"

# Match the saved selection / active cell
$ws.Range("B1").Select()
